# Common: So, liquid creation is backed by transactions from pricelist
#
# Adds a new "default" tariff row (code JX-BMD9-GYJXO7) to the "tariffs"
# sheet, right below the existing "default" row, and moves the selection
# to C4 to mirror the state the workbook was saved in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "default"
$ws.Range("B3").Value = "market.tariff.default"
$ws.Range("C3").Value = "JX-BMD9-GYJXO7"

$ws.Range("C4").Select()
